$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value, whether it is a numeric-looking string that
# must be forced to Text to avoid Excel auto-converting it to a number (which would
# drop formatting like trailing zeros, e.g. "1.00" -> 1).
$updates = @(
    @{ Cell = 'D2'; Value = '93.293.70'; Numeric = $false }
    @{ Cell = 'E2'; Value = '  -5.36%  '; Numeric = $false }
    @{ Cell = 'D3'; Value = '3.343.68'; Numeric = $false }
    @{ Cell = 'E3'; Value = '  -3.83%  '; Numeric = $false }
    @{ Cell = 'D4'; Value = '1.00'; Numeric = $true }
    @{ Cell = 'E4'; Value = '  +0.10%  '; Numeric = $false }
    @{ Cell = 'D5'; Value = '232.42'; Numeric = $true }
    @{ Cell = 'E5'; Value = '  -8.86%  '; Numeric = $false }
    @{ Cell = 'D6'; Value = '627.75'; Numeric = $true }
    @{ Cell = 'E6'; Value = '  -6.04%  '; Numeric = $false }
    @{ Cell = 'D7'; Value = '1.36'; Numeric = $true }
    @{ Cell = 'E7'; Value = '  -9.48%  '; Numeric = $false }
    @{ Cell = 'D8'; Value = '0.387'; Numeric = $true }
    @{ Cell = 'E8'; Value = '  -10.12%  '; Numeric = $false }
    @{ Cell = 'E9'; Value = '  +0.14%  '; Numeric = $false }
    @{ Cell = 'D10'; Value = '0.936'; Numeric = $true }
    @{ Cell = 'E10'; Value = '  -11.51%  '; Numeric = $false }
    @{ Cell = 'D11'; Value = '3.342.01'; Numeric = $false }
    @{ Cell = 'E11'; Value = '  -3.81%  '; Numeric = $false }
    @{ Cell = 'E12'; Value = '  -7.59%  '; Numeric = $false }
    @{ Cell = 'D13'; Value = '40.22'; Numeric = $true }
    @{ Cell = 'E13'; Value = '  -12.82%  '; Numeric = $false }
    @{ Cell = 'D14'; Value = '5.98'; Numeric = $true }
    @{ Cell = 'E14'; Value = '  -3.55%  '; Numeric = $false }
    @{ Cell = 'D15'; Value = '93.195.47'; Numeric = $false }
    @{ Cell = 'E15'; Value = '  -5.31%  '; Numeric = $false }
    @{ Cell = 'D16'; Value = '3.977.56'; Numeric = $false }
    @{ Cell = 'E16'; Value = '  -3.40%  '; Numeric = $false }
    @{ Cell = 'D17'; Value = '0.0000243'; Numeric = $true }
    @{ Cell = 'E17'; Value = '  -6.69%  '; Numeric = $false }
    @{ Cell = 'D18'; Value = '7.97'; Numeric = $true }
    @{ Cell = 'E18'; Value = '  -11.87%  '; Numeric = $false }
    @{ Cell = 'D19'; Value = '3.355.19'; Numeric = $false }
    @{ Cell = 'E19'; Value = '  -2.51%  '; Numeric = $false }
    @{ Cell = 'D20'; Value = '16.83'; Numeric = $true }
    @{ Cell = 'E20'; Value = '  -9.78%  '; Numeric = $false }
    @{ Cell = 'D21'; Value = '10.88'; Numeric = $true }
    @{ Cell = 'E21'; Value = '  -7.62%  '; Numeric = $false }
    @{ Cell = 'D22'; Value = '493.05'; Numeric = $true }
    @{ Cell = 'E22'; Value = '  -5.24%  '; Numeric = $false }
    @{ Cell = 'D23'; Value = '0.451'; Numeric = $true }
    @{ Cell = 'E23'; Value = '  -15.89%  '; Numeric = $false }
    @{ Cell = 'D24'; Value = '3.13'; Numeric = $true }
    @{ Cell = 'E24'; Value = '  -9.63%  '; Numeric = $false }
    @{ Cell = 'D25'; Value = '0.0000185'; Numeric = $true }
    @{ Cell = 'E25'; Value = '  -9.02%  '; Numeric = $false }
    @{ Cell = 'D26'; Value = '6.17'; Numeric = $true }
    @{ Cell = 'E26'; Value = '  -9.85%  '; Numeric = $false }
    @{ Cell = 'D27'; Value = '89.78'; Numeric = $true }
    @{ Cell = 'E27'; Value = '  -8.39%  '; Numeric = $false }
    @{ Cell = 'D28'; Value = '3.532.17'; Numeric = $false }
    @{ Cell = 'E28'; Value = '  -3.35%  '; Numeric = $false }
    @{ Cell = 'D29'; Value = '11.45'; Numeric = $true }
    @{ Cell = 'E29'; Value = '  -9.54%  '; Numeric = $false }
    @{ Cell = 'D30'; Value = '11.31'; Numeric = $true }
    @{ Cell = 'E30'; Value = '  -9.25%  '; Numeric = $false }
    @{ Cell = 'D31'; Value = '1.00'; Numeric = $true }
    @{ Cell = 'E31'; Value = '  +0.09%  '; Numeric = $false }
    @{ Cell = 'D32'; Value = '2.64'; Numeric = $true }
    @{ Cell = 'E32'; Value = '  -7.61%  '; Numeric = $false }
    @{ Cell = 'E33'; Value = '  -9.96%  '; Numeric = $false }
    @{ Cell = 'D34'; Value = '0.998'; Numeric = $true }
    @{ Cell = 'E34'; Value = '  +0.06%  '; Numeric = $false }
    @{ Cell = 'E35'; Value = '  -10.24%  '; Numeric = $false }
    @{ Cell = 'D36'; Value = '28.60'; Numeric = $true }
    @{ Cell = 'E36'; Value = '  -5.22%  '; Numeric = $false }
    @{ Cell = 'D37'; Value = '0.525'; Numeric = $true }
    @{ Cell = 'E37'; Value = '  -10.16%  '; Numeric = $false }
    @{ Cell = 'D38'; Value = '7.46'; Numeric = $true }
    @{ Cell = 'E38'; Value = '  -7.32%  '; Numeric = $false }
    @{ Cell = 'D39'; Value = '522.06'; Numeric = $true }
    @{ Cell = 'E39'; Value = '  -1.73%  '; Numeric = $false }
    @{ Cell = 'E40'; Value = '  +0.09%  '; Numeric = $false }
    @{ Cell = 'D41'; Value = '1.39'; Numeric = $true }
    @{ Cell = 'E41'; Value = '  -8.86%  '; Numeric = $false }
    @{ Cell = 'D42'; Value = '0.148'; Numeric = $true }
    @{ Cell = 'E42'; Value = '  -5.31%  '; Numeric = $false }
    @{ Cell = 'D43'; Value = '0.873'; Numeric = $true }
    @{ Cell = 'E43'; Value = '  -2.11%  '; Numeric = $false }
    @{ Cell = 'D44'; Value = '24.03'; Numeric = $true }
    @{ Cell = 'E44'; Value = '  -1.69%  '; Numeric = $false }
    @{ Cell = 'D45'; Value = '1.68'; Numeric = $true }
    @{ Cell = 'E45'; Value = '  -5.90%  '; Numeric = $false }
    @{ Cell = 'B46'; Value = 'Filecoin'; Numeric = $false }
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; Numeric = $false }
    @{ Cell = 'D46'; Value = '5.48'; Numeric = $true }
    @{ Cell = 'E46'; Value = '  -5.09%  '; Numeric = $false }
    @{ Cell = 'B47'; Value = 'MantraDAO'; Numeric = $false }
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'; Numeric = $false }
    @{ Cell = 'D47'; Value = '3.56'; Numeric = $true }
    @{ Cell = 'E47'; Value = '  -3.07%  '; Numeric = $false }
    @{ Cell = 'D48'; Value = '2.15'; Numeric = $true }
    @{ Cell = 'E48'; Value = '  -4.96%  '; Numeric = $false }
    @{ Cell = 'D49'; Value = '0.0394'; Numeric = $true }
    @{ Cell = 'E49'; Value = '  -9.66%  '; Numeric = $false }
    @{ Cell = 'D50'; Value = '52.36'; Numeric = $true }
    @{ Cell = 'E50'; Value = '  -5.94%  '; Numeric = $false }
    @{ Cell = 'D51'; Value = '3.10'; Numeric = $true }
    @{ Cell = 'E51'; Value = '  -4.84%  '; Numeric = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.Numeric) {
        # Pre-format as text so the numeric-looking string is stored verbatim
        # (preserving trailing zeros / exact decimal text) instead of becoming a
        # floating-point number, then drop back to the default "Normal" style so
        # no stray number-format is left behind on the cell.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
